$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.048.53"
$ws.Range("E2").Value = "  +7.72%  "
$ws.Range("D3").Value = "3.368.90"
$ws.Range("E3").Value = "  +6.17%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.80"
$ws.Range("E5").Value = "  +3.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "639.80"
$ws.Range("E6").Value = "  +2.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.415"
$ws.Range("E7").Value = "  +47.27%  "
$ws.Range("E8").Value = "  +13.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "3.367.38"
$ws.Range("E10").Value = "  +6.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.613"
$ws.Range("E11").Value = "  +4.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000280"
$ws.Range("E12").Value = "  +10.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.92"
$ws.Range("E13").Value = "  +14.19%  "
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").Value = "3.982.40"
$ws.Range("E15").Value = "  +6.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "87.956.99"
$ws.Range("E17").Value = "  +8.03%  "
$ws.Range("D18").Value = "3.377.46"
$ws.Range("E18").Value = "  +6.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.81"
$ws.Range("E19").Value = "  +6.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.55"
$ws.Range("E20").Value = "  +7.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "451.71"
$ws.Range("E21").Value = "  +4.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.01"
$ws.Range("E22").Value = "  -5.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.50"
$ws.Range("E23").Value = "  +8.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.39"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("E25").Value = "  +3.93%  "
$ws.Range("E26").Value = "  +14.35%  "
$ws.Range("D27").Value = "3.563.03"
$ws.Range("E27").Value = "  +7.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "80.18"
$ws.Range("E28").Value = "  +4.85%  "
$ws.Range("E29").Value = "  +15.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.184"
$ws.Range("E31").Value = "  +33.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.24"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "569.18"
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.52"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.07"
$ws.Range("E36").Value = "  +4.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.22"
$ws.Range("E37").Value = "  +17.61%  "
$ws.Range("E38").Value = "  -8.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.64"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.429"
$ws.Range("E40").Value = "  +5.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.85"
$ws.Range("E41").Value = "  +5.27%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.06"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "157.30"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "185.87"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.40"
$ws.Range("E48").Value = "  +4.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.26"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("E50").Value = "  +6.12%  "
$ws.Range("E51").Value = "  +5.19%  "
